$wb = $excel.ActiveWorkbook

# Sheet 1: sum_response_time_Results - update row 4 (Use Case 3)
$ws1 = $wb.Worksheets.Item("sum_response_time_Results")
$ws1.Range("C4").Value = 15.16
$ws1.Range("D4").Value = 5.296640444659237
$ws1.Range("G4").Value = [double]"8.154447662214437E-06"
$ws1.Range("K4").Value = [double]"3.929260964413422E-12"

# Sheet 2: total_data_transferred_Results - update row 4 (Use Case 3)
$ws2 = $wb.Worksheets.Item("total_data_transferred_Results")
$ws2.Range("C4").Value = 3560
